$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 111112890
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 250001500
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 250001500
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -250001850

# Row 112
$ws.Range("H112").Value = 2288.5293
$ws.Range("J112").Value = 2631.3794
$ws.Range("L112").Value = 7894.138199999999
$ws.Range("N112").Value = -10110.1382

# Row 137
$ws.Range("H137").Value = 3948728.5
$ws.Range("I137").Value = 1924318.4
$ws.Range("J137").Value = 8334950
$ws.Range("K137").Value = 5772955.199999999
$ws.Range("L137").Value = 25004850
$ws.Range("M137").Value = -5770405.199999999
$ws.Range("N137").Value = -25009950

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 1034.1936
$ws.Range("I74").Value = 1046.1482
$ws.Range("K74").Value = 1046.1482
$ws.Range("M74").Value = -172.1482000000001

# Row 77
$ws.Range("H77").Value = 1034.1936
$ws.Range("I77").Value = 1046.1482
$ws.Range("K77").Value = 5230.741
$ws.Range("M77").Value = -862.741

# Row 102
$ws.Range("H102").Value = 83335420
$ws.Range("I102").Value = 2260.9
$ws.Range("J102").Value = 500001250
$ws.Range("K102").Value = 2260.9
$ws.Range("L102").Value = 500001250
$ws.Range("M102").Value = -638.9000000000001
$ws.Range("N102").Value = -500004494

# Row 110
$ws.Range("H110").Value = 5955.0527
$ws.Range("I110").Value = 5598.8667
$ws.Range("J110").Value = 7290.75
$ws.Range("K110").Value = 5598.8667
$ws.Range("L110").Value = 7290.75
$ws.Range("M110").Value = -3553.8667
$ws.Range("N110").Value = -11380.75

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 865.6727
$ws.Range("I58").Value = 896.8837
$ws.Range("J58").Value = 753.8333
$ws.Range("K58").Value = 896.8837
$ws.Range("L58").Value = 753.8333
$ws.Range("M58").Value = -693.8837
$ws.Range("N58").Value = -1159.8333

# Row 136
$ws.Range("H136").Value = 865.6727
$ws.Range("I136").Value = 896.8837
$ws.Range("J136").Value = 753.8333
$ws.Range("K136").Value = 2690.6511
$ws.Range("L136").Value = 2261.4999
$ws.Range("M136").Value = -140.6511
$ws.Range("N136").Value = -7361.4999

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 834153.8
$ws.Range("I5").Value = 780.5714
$ws.Range("J5").Value = 6667766.5
$ws.Range("K5").Value = 2341.7142
$ws.Range("L5").Value = 20003299.5
$ws.Range("M5").Value = -2229.7142
$ws.Range("N5").Value = -20003523.5

# Row 63
$ws.Range("H63").Value = 4232.3335
$ws.Range("I63").Value = 2514.6667
$ws.Range("J63").Value = 5950
$ws.Range("K63").Value = 7544.000100000001
$ws.Range("L63").Value = 17850
$ws.Range("M63").Value = -6795.000100000001
$ws.Range("N63").Value = -19348

# Row 64
$ws.Range("H64").Value = 6616.364
$ws.Range("I64").Value = 860.5
$ws.Range("J64").Value = 7895.4443
$ws.Range("K64").Value = 2581.5
$ws.Range("L64").Value = 23686.3329
$ws.Range("M64").Value = -2311.5
$ws.Range("N64").Value = -24226.3329

# Row 66
$ws.Range("H66").Value = 4232.3335
$ws.Range("I66").Value = 2514.6667
$ws.Range("J66").Value = 5950
$ws.Range("K66").Value = 22632.0003
$ws.Range("L66").Value = 53550
$ws.Range("M66").Value = -18888.0003
$ws.Range("N66").Value = -61038

# Row 67
$ws.Range("H67").Value = 6616.364
$ws.Range("I67").Value = 860.5
$ws.Range("J67").Value = 7895.4443
$ws.Range("K67").Value = 2581.5
$ws.Range("L67").Value = 23686.3329
$ws.Range("M67").Value = -1645.5
$ws.Range("N67").Value = -25558.3329

# Row 86
$ws.Range("H86").Value = 1390.3
$ws.Range("I86").Value = 1066.6666
$ws.Range("J86").Value = 1529
$ws.Range("K86").Value = 3199.9998
$ws.Range("L86").Value = 4587
$ws.Range("M86").Value = -2013.9998
$ws.Range("N86").Value = -6959

# Row 89
$ws.Range("H89").Value = 1390.3
$ws.Range("I89").Value = 1066.6666
$ws.Range("J89").Value = 1529
$ws.Range("K89").Value = 9599.999400000001
$ws.Range("L89").Value = 13761
$ws.Range("M89").Value = -3671.999400000001
$ws.Range("N89").Value = -25617

# Row 114
$ws.Range("H114").Value = 1919.8695
$ws.Range("I114").Value = 174.38461
$ws.Range("J114").Value = 4189
$ws.Range("K114").Value = 523.15383
$ws.Range("L114").Value = 12567
$ws.Range("M114").Value = 2730.84617
$ws.Range("N114").Value = -19075

# Row 122
$ws.Range("H122").Value = 45064.652
$ws.Range("I122").Value = 318.7
$ws.Range("J122").Value = 50328.883
$ws.Range("K122").Value = 2868.3
$ws.Range("L122").Value = 452959.947
$ws.Range("M122").Value = -418.2999999999997
$ws.Range("N122").Value = -457859.947

# Row 132
$ws.Range("H132").Value = 732511.7
$ws.Range("I132").Value = 1013093.06
$ws.Range("K132").Value = 9117837.540000001
$ws.Range("M132").Value = -9115307.540000001

# Row 135
$ws.Range("H135").Value = 834153.8
$ws.Range("I135").Value = 780.5714
$ws.Range("J135").Value = 6667766.5
$ws.Range("K135").Value = 7025.1426
$ws.Range("L135").Value = 60009898.5
$ws.Range("M135").Value = -4490.1426
$ws.Range("N135").Value = -60014968.5

# Row 137
$ws.Range("H137").Value = 2206.6365
$ws.Range("I137").Value = 1905.7142
$ws.Range("J137").Value = 2733.25
$ws.Range("K137").Value = 5717.142599999999
$ws.Range("L137").Value = 8199.75
$ws.Range("M137").Value = -617.1425999999992
$ws.Range("N137").Value = -18399.75

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2282.647
$ws.Range("J61").Value = 3451.25
$ws.Range("L61").Value = 3451.25
$ws.Range("N61").Value = -3855.25

# Row 113
$ws.Range("H113").Value = 2282.647
$ws.Range("J113").Value = 3451.25
$ws.Range("L113").Value = 3451.25
$ws.Range("N113").Value = -7791.25

# Row 136
$ws.Range("H136").Value = 2488.0557
$ws.Range("I136").Value = 2056.0715
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 6168.2145
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -3618.2145
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1893
$ws.Range("I126").Value = 1530.9131
$ws.Range("J126").Value = 3975
$ws.Range("K126").Value = 4592.7393
$ws.Range("L126").Value = 11925
$ws.Range("M126").Value = -2122.7393
$ws.Range("N126").Value = -16865

# Row 132
$ws.Range("H132").Value = 2985.5
$ws.Range("J132").Value = 2577.6
$ws.Range("L132").Value = 7732.799999999999
$ws.Range("N132").Value = -12792.8
